# Atualização automática: 2025-08-20 09:01:05
#
# Rows 7-11 hold one detection-image "batch" each; the batch boundary
# moved by one row (old row 7 -> new row 11, old rows 8-11 -> new rows
# 7-10). Columns A and D-J carry the per-batch data; B (Class) and C
# (First_Detection_Date) are identical for every row in the block, so
# they are left untouched.
#
# Cell-to-cell Copy/PasteSpecial (values) is used instead of reading
# .Value into a variable and writing it back: some of the bounding-box
# strings in column I (e.g. "702,633,740,690") parse as a
# thousands-grouped number, and a plain Value round-trip would silently
# turn them into numerics. Copy/PasteSpecial preserves the original
# text cell type. A scratch row well below the used range is used to
# hold the wrap-around value (old row 7) while rows 8-11 shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$scratchRow = 20

# (source row, destination row) pairs, in the order they must be applied
# so that nothing is overwritten before it has been copied out:
#   7 -> scratch (save the wrap-around batch)
#   8 -> 7, 9 -> 8, 10 -> 9, 11 -> 10 (shift up)
#   scratch -> 11 (restore the wrap-around batch)
$moves = @(
    @(7, $scratchRow),
    @(8, 7),
    @(9, 8),
    @(10, 9),
    @(11, 10),
    @($scratchRow, 11)
)

foreach ($move in $moves) {
    $srcRow = $move[0]
    $dstRow = $move[1]

    $ws.Range("A$srcRow").Copy() | Out-Null
    $ws.Range("A$dstRow").PasteSpecial($xlPasteValues) | Out-Null

    $ws.Range("D${srcRow}:J${srcRow}").Copy() | Out-Null
    $ws.Range("D${dstRow}:J${dstRow}").PasteSpecial($xlPasteValues) | Out-Null
}

$ws.Range("A${scratchRow}:J${scratchRow}").ClearContents() | Out-Null
$excel.CutCopyMode = $false

# Row 18: refreshed detection image / bounding box / confidence.
# D18/I18 are plain (non-numeric-looking) text, so a direct .Value
# assignment keeps its inline-string type. J18's new value ("0.76") does
# look numeric, so a direct assignment would silently store it as a
# number instead of text (unlike the original "0.75" inline string) -
# copy it instead from J4, which already holds the literal text "0.76",
# the same way the row 7-11 shift above avoids that pitfall.
$ws.Range("D18").Value = "image_20250808221835_ppp0.jpg"
$ws.Range("I18").Value = "1182,405,1231,455"

$ws.Range("J4").Copy() | Out-Null
$ws.Range("J18").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
